$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update current page bookmark for "Researching Information Systems and Computing" (row 11)
$ws.Range("C11").Value = 175

# Recalculate so the SUM formula in E7 reflects the new total page count
$excel.Calculate()

# Update the active selection to match the author's last cursor position
$ws.Range("E15").Select()
